# Updated cryptos list on Tue Jul 25 10:29:48 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns for rows 2-51 with
# the latest scraped values. Values are kept as literal text (matching the
# existing inline-string cell contents, e.g. "1.0000", "29.187.59", the
# padded "  -0.66%  " percentage strings, etc.), so the target cells are
# forced to Text number format first -- otherwise Excel would silently
# reinterpret numeric-looking strings (like "1.0000" -> 1, or "0.9996")
# as numbers and normalize/round them away from the source formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.187.59' },
    @{ Cell = 'E2'; Value = '  -0.66%  ' },
    @{ Cell = 'D3'; Value = '1.856.77' },
    @{ Cell = 'E3'; Value = '  +0.17%  ' },
    @{ Cell = 'D4'; Value = '0.9996' },
    @{ Cell = 'E4'; Value = '  -0.23%  ' },
    @{ Cell = 'D5'; Value = '237.97' },
    @{ Cell = 'E5'; Value = '  -0.41%  ' },
    @{ Cell = 'D6'; Value = '0.6904' },
    @{ Cell = 'E6'; Value = '  -1.65%  ' },
    @{ Cell = 'D7'; Value = '1.0000' },
    @{ Cell = 'E7'; Value = '  -0.17%  ' },
    @{ Cell = 'D8'; Value = '0.07777' },
    @{ Cell = 'E8'; Value = '  +5.30%  ' },
    @{ Cell = 'E9'; Value = '  -0.67%  ' },
    @{ Cell = 'E10'; Value = '  -1.91%  ' },
    @{ Cell = 'D11'; Value = '0.08072' },
    @{ Cell = 'E11'; Value = '  -0.55%  ' },
    @{ Cell = 'D12'; Value = '1.864.44' },
    @{ Cell = 'E12'; Value = '  -0.14%  ' },
    @{ Cell = 'E13'; Value = '  -0.63%  ' },
    @{ Cell = 'D14'; Value = '5.185' },
    @{ Cell = 'E14'; Value = '  -0.54%  ' },
    @{ Cell = 'D15'; Value = '89.41' },
    @{ Cell = 'E15'; Value = '  -0.37%  ' },
    @{ Cell = 'D16'; Value = '29.196.40' },
    @{ Cell = 'E16'; Value = '  -1.97%  ' },
    @{ Cell = 'D17'; Value = '5.738' },
    @{ Cell = 'E17'; Value = '  -2.80%  ' },
    @{ Cell = 'D18'; Value = '0.000007811' },
    @{ Cell = 'E18'; Value = '  +0.90%  ' },
    @{ Cell = 'E19'; Value = '  +1.00%  ' },
    @{ Cell = 'D20'; Value = '234.94' },
    @{ Cell = 'E20'; Value = '  -3.15%  ' },
    @{ Cell = 'D21'; Value = '0.9998' },
    @{ Cell = 'E21'; Value = '  -0.40%  ' },
    @{ Cell = 'D22'; Value = '2.113.56' },
    @{ Cell = 'E22'; Value = '  -3.00%  ' },
    @{ Cell = 'D23'; Value = '0.9995' },
    @{ Cell = 'E23'; Value = '  -0.32%  ' },
    @{ Cell = 'D24'; Value = '7.477' },
    @{ Cell = 'E24'; Value = '  -2.09%  ' },
    @{ Cell = 'D25'; Value = '161.88' },
    @{ Cell = 'E25'; Value = '  +0.31%  ' },
    @{ Cell = 'D26'; Value = '8.969' },
    @{ Cell = 'E26'; Value = '  -0.82%  ' },
    @{ Cell = 'D27'; Value = '0.1426' },
    @{ Cell = 'E27'; Value = '  -3.96%  ' },
    @{ Cell = 'D28'; Value = '18.05' },
    @{ Cell = 'E28'; Value = '  -0.24%  ' },
    @{ Cell = 'E29'; Value = '  +0.45%  ' },
    @{ Cell = 'D30'; Value = '1.399' },
    @{ Cell = 'E30'; Value = '  +0.68%  ' },
    @{ Cell = 'D31'; Value = '4.529' },
    @{ Cell = 'E31'; Value = '  +2.59%  ' },
    @{ Cell = 'D32'; Value = '1.482' },
    @{ Cell = 'E32'; Value = '  -1.64%  ' },
    @{ Cell = 'D33'; Value = '4.011' },
    @{ Cell = 'E33'; Value = '  -1.41%  ' },
    @{ Cell = 'D34'; Value = '0.05201' },
    @{ Cell = 'E34'; Value = '  -2.23%  ' },
    @{ Cell = 'D35'; Value = '1.185' },
    @{ Cell = 'E35'; Value = '  -1.32%  ' },
    @{ Cell = 'D36'; Value = '0.7033' },
    @{ Cell = 'E36'; Value = '  -2.91%  ' },
    @{ Cell = 'D37'; Value = '1.026' },
    @{ Cell = 'E37'; Value = '  +1.79%  ' },
    @{ Cell = 'D38'; Value = '2.671' },
    @{ Cell = 'E38'; Value = '  -0.42%  ' },
    @{ Cell = 'D39'; Value = '0.01846' },
    @{ Cell = 'E39'; Value = '  -1.28%  ' },
    @{ Cell = 'D40'; Value = '2.677' },
    @{ Cell = 'E40'; Value = '  -1.66%  ' },
    @{ Cell = 'D41'; Value = '0.9298' },
    @{ Cell = 'E41'; Value = '  +6.27%  ' },
    @{ Cell = 'D42'; Value = '1.092.06' },
    @{ Cell = 'E42'; Value = '  +5.84%  ' },
    @{ Cell = 'D43'; Value = '5.995' },
    @{ Cell = 'E43'; Value = '  +1.00%  ' },
    @{ Cell = 'D44'; Value = '0.4288' },
    @{ Cell = 'E44'; Value = '  -0.85%  ' },
    @{ Cell = 'D45'; Value = '70.52' },
    @{ Cell = 'E45'; Value = '  +1.08%  ' },
    @{ Cell = 'D46'; Value = '0.9998' },
    @{ Cell = 'E46'; Value = '  -0.20%  ' },
    @{ Cell = 'D47'; Value = '102.70' },
    @{ Cell = 'E47'; Value = '  +0.14%  ' },
    @{ Cell = 'D48'; Value = '1.796' },
    @{ Cell = 'E48'; Value = '  +2.01%  ' },
    @{ Cell = 'D49'; Value = '2.008.42' },
    @{ Cell = 'E49'; Value = '  -2.87%  ' },
    @{ Cell = 'D50'; Value = '9.158' },
    @{ Cell = 'E50'; Value = '  -0.16%  ' },
    @{ Cell = 'D51'; Value = '7.003' },
    @{ Cell = 'E51'; Value = '  -3.71%  ' },
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}

Write-Output ("Updated {0} cells" -f $updates.Count)
